$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Cells.Item(8, 1).Value = 46045
$ws.Cells.Item(8, 4).Value = 157.83000000000001
$ws.Cells.Item(8, 5).Value = 150.65
$ws.Cells.Item(8, 6).Value = 160.65
$ws.Cells.Item(8, 7).Value = 150.66
# Row 9
$ws.Cells.Item(9, 1).Value = 46045
$ws.Cells.Item(9, 4).Value = 157.83000000000001
$ws.Cells.Item(9, 5).Value = 150.65
$ws.Cells.Item(9, 6).Value = 160.65
$ws.Cells.Item(9, 7).Value = 150.66
# Row 10
$ws.Cells.Item(10, 1).Value = 46045
$ws.Cells.Item(10, 4).Value = 158.59
$ws.Cells.Item(10, 5).Value = 152.4
$ws.Cells.Item(10, 6).Value = 162.4
$ws.Cells.Item(10, 7).Value = 152.80000000000001
# Row 11
$ws.Cells.Item(11, 1).Value = 46044
$ws.Cells.Item(11, 4).Value = 157.58000000000001
$ws.Cells.Item(11, 5).Value = 150.74
$ws.Cells.Item(11, 6).Value = 160.74
$ws.Cells.Item(11, 7).Value = 150.75
# Row 12
$ws.Cells.Item(12, 1).Value = 46044
$ws.Cells.Item(12, 4).Value = 157.58000000000001
$ws.Cells.Item(12, 5).Value = 150.74
$ws.Cells.Item(12, 6).Value = 160.74
$ws.Cells.Item(12, 7).Value = 150.75
# Row 13
$ws.Cells.Item(13, 1).Value = 46044
$ws.Cells.Item(13, 4).Value = 158.36000000000001
$ws.Cells.Item(13, 5).Value = 152.51
$ws.Cells.Item(13, 6).Value = 162.51
$ws.Cells.Item(13, 7).Value = 152.91
# Row 17
$ws.Cells.Item(17, 1).Value = 46045
$ws.Cells.Item(17, 4).Value = 161.97
$ws.Cells.Item(17, 5).Value = 154.88
$ws.Cells.Item(17, 6).Value = 164.88
# Row 18
$ws.Cells.Item(18, 1).Value = 46044
$ws.Cells.Item(18, 4).Value = 161.74
$ws.Cells.Item(18, 5).Value = 155
$ws.Cells.Item(18, 6).Value = 165
# Row 22
$ws.Cells.Item(22, 1).Value = 46045
$ws.Cells.Item(22, 4).Value = 158.91
$ws.Cells.Item(22, 5).Value = 152.22
$ws.Cells.Item(22, 6).Value = 161.82
$ws.Cells.Item(22, 7).Value = 153.29
# Row 23
$ws.Cells.Item(23, 1).Value = 46045
$ws.Cells.Item(23, 4).Value = 163.37
$ws.Cells.Item(23, 5).Value = 157.72999999999999
$ws.Cells.Item(23, 6).Value = 167.73
# Row 24
$ws.Cells.Item(24, 1).Value = 46045
$ws.Cells.Item(24, 4).Value = 163.5
$ws.Cells.Item(24, 5).Value = 158.46
$ws.Cells.Item(24, 6).Value = 168.46
# Row 25
$ws.Cells.Item(25, 1).Value = 46045
$ws.Cells.Item(25, 4).Value = 163.49
$ws.Cells.Item(25, 5).Value = 158
$ws.Cells.Item(25, 6).Value = 168
$ws.Cells.Item(25, 7).Value = 158.13
# Row 26
$ws.Cells.Item(26, 1).Value = 46045
$ws.Cells.Item(26, 4).Value = 163.06
$ws.Cells.Item(26, 5).Value = 159.62
$ws.Cells.Item(26, 6).Value = 169.62
# Row 27
$ws.Cells.Item(27, 1).Value = 46044
$ws.Cells.Item(27, 4).Value = 158.65
$ws.Cells.Item(27, 5).Value = 152.31
$ws.Cells.Item(27, 6).Value = 161.91
$ws.Cells.Item(27, 7).Value = 153.38
# Row 28
$ws.Cells.Item(28, 1).Value = 46044
$ws.Cells.Item(28, 4).Value = 163.13999999999999
$ws.Cells.Item(28, 5).Value = 157.84
$ws.Cells.Item(28, 6).Value = 167.84
# Row 29
$ws.Cells.Item(29, 1).Value = 46044
$ws.Cells.Item(29, 4).Value = 163.27000000000001
$ws.Cells.Item(29, 5).Value = 158.56
$ws.Cells.Item(29, 6).Value = 168.56
# Row 30
$ws.Cells.Item(30, 1).Value = 46044
$ws.Cells.Item(30, 4).Value = 163.25
$ws.Cells.Item(30, 5).Value = 158.11000000000001
$ws.Cells.Item(30, 6).Value = 168.11
$ws.Cells.Item(30, 7).Value = 158.22999999999999
# Row 31
$ws.Cells.Item(31, 1).Value = 46044
$ws.Cells.Item(31, 4).Value = 162.83000000000001
$ws.Cells.Item(31, 5).Value = 159.72
$ws.Cells.Item(31, 6).Value = 169.72
# Row 35
$ws.Cells.Item(35, 1).Value = 46045
$ws.Cells.Item(35, 4).Value = 157.47999999999999
$ws.Cells.Item(35, 5).Value = 149.38999999999999
$ws.Cells.Item(35, 6).Value = 158.38999999999999
# Row 36
$ws.Cells.Item(36, 1).Value = 46044
$ws.Cells.Item(36, 4).Value = 157.24
$ws.Cells.Item(36, 5).Value = 149.5
$ws.Cells.Item(36, 6).Value = 158.5
# Row 40
$ws.Cells.Item(40, 1).Value = 46045
$ws.Cells.Item(40, 4).Value = 163.05000000000001
$ws.Cells.Item(40, 5).Value = 157.75
$ws.Cells.Item(40, 6).Value = 167.75
# Row 41
$ws.Cells.Item(41, 1).Value = 46045
$ws.Cells.Item(41, 4).Value = 162.76
$ws.Cells.Item(41, 5).Value = 158.16999999999999
$ws.Cells.Item(41, 6).Value = 168.17
# Row 42
$ws.Cells.Item(42, 1).Value = 46044
$ws.Cells.Item(42, 4).Value = 162.80000000000001
$ws.Cells.Item(42, 5).Value = 157.80000000000001
$ws.Cells.Item(42, 6).Value = 167.8
# Row 43
$ws.Cells.Item(43, 1).Value = 46044
$ws.Cells.Item(43, 4).Value = 162.52000000000001
$ws.Cells.Item(43, 5).Value = 158.22
$ws.Cells.Item(43, 6).Value = 168.22
# Row 47
$ws.Cells.Item(47, 1).Value = 46045
$ws.Cells.Item(47, 4).Value = 157.1
$ws.Cells.Item(47, 5).Value = 151.18
$ws.Cells.Item(47, 6).Value = 161.18
# Row 48
$ws.Cells.Item(48, 1).Value = 46045
$ws.Cells.Item(48, 4).Value = 156.68
$ws.Cells.Item(48, 5).Value = 151.08000000000001
$ws.Cells.Item(48, 6).Value = 161.08000000000001
# Row 49
$ws.Cells.Item(49, 1).Value = 46044
$ws.Cells.Item(49, 4).Value = 156.76
$ws.Cells.Item(49, 5).Value = 151.16
$ws.Cells.Item(49, 6).Value = 161.16
# Row 50
$ws.Cells.Item(50, 1).Value = 46044
$ws.Cells.Item(50, 4).Value = 156.35
$ws.Cells.Item(50, 5).Value = 151.07
$ws.Cells.Item(50, 6).Value = 161.07
# Row 54
$ws.Cells.Item(54, 1).Value = 46045
$ws.Cells.Item(54, 4).Value = 172.03
$ws.Cells.Item(54, 5).Value = 165.81
$ws.Cells.Item(54, 6).Value = 175.81
# Row 55
$ws.Cells.Item(55, 1).Value = 46045
$ws.Cells.Item(55, 4).Value = 164.86
$ws.Cells.Item(55, 5).Value = 163.58000000000001
$ws.Cells.Item(55, 6).Value = 173.58
# Row 56
$ws.Cells.Item(56, 1).Value = 46045
$ws.Cells.Item(56, 4).Value = 161.55000000000001
# Row 57
$ws.Cells.Item(57, 1).Value = 46045
$ws.Cells.Item(57, 4).Value = 161.86000000000001
$ws.Cells.Item(57, 5).Value = 158
# Row 58
$ws.Cells.Item(58, 1).Value = 46045
$ws.Cells.Item(58, 4).Value = 157.63
$ws.Cells.Item(58, 5).Value = 153.9
$ws.Cells.Item(58, 6).Value = 163.9
# Row 59
$ws.Cells.Item(59, 1).Value = 46045
$ws.Cells.Item(59, 4).Value = 164.28
$ws.Cells.Item(59, 5).Value = 163.79
# Row 60
$ws.Cells.Item(60, 1).Value = 46044
$ws.Cells.Item(60, 4).Value = 171.8
$ws.Cells.Item(60, 5).Value = 165.87
$ws.Cells.Item(60, 6).Value = 175.87
# Row 61
$ws.Cells.Item(61, 1).Value = 46044
$ws.Cells.Item(61, 4).Value = 164.61
$ws.Cells.Item(61, 5).Value = 163.69999999999999
$ws.Cells.Item(61, 6).Value = 173.7
# Row 62
$ws.Cells.Item(62, 1).Value = 46044
$ws.Cells.Item(62, 4).Value = 161.32
# Row 63
$ws.Cells.Item(63, 1).Value = 46044
$ws.Cells.Item(63, 4).Value = 161.63999999999999
$ws.Cells.Item(63, 5).Value = 158.12
# Row 64
$ws.Cells.Item(64, 1).Value = 46044
$ws.Cells.Item(64, 4).Value = 157.41
$ws.Cells.Item(64, 5).Value = 154.02000000000001
$ws.Cells.Item(64, 6).Value = 164.02
# Row 65
$ws.Cells.Item(65, 1).Value = 46044
$ws.Cells.Item(65, 4).Value = 164.05
$ws.Cells.Item(65, 5).Value = 163.87
